# v0.2 - opens and scrapes needed metadata to excel
#
# The sheet previously contained a header ("links") in A1 followed by 10
# scraped Play Store URLs in A2:A11. Several of those scraped rows
# (originally A5:A10) are removed, leaving only the first three scraped
# links plus the last one (originally A11, "...com.matteljv.uno"), which
# shifts up to become the new A5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the no-longer-needed rows; remaining rows below shift up
# automatically, and the shared string table is compacted on save.
$ws.Rows("5:10").Delete()

# Match the active selection left behind in the saved workbook.
$ws.Range("G5").Select() | Out-Null
